# Update MOSIP_Partner Management Requirements.xlsx
# Adds a new "Estimates" worksheet (after Sheet1) containing a single
# labeled estimate row, and makes that new sheet the active tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1, rename it "Estimates"
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Estimates"

# Column widths to match the target layout (OOXML stored widths: B=56, C~10.54 chars)
$newSheet.Columns.Item(2).ColumnWidth = 55.166666666666664
$newSheet.Columns.Item(3).ColumnWidth = 9.65

# Row height for row 2 (tall enough to show the wrapped label in 3 lines)
$newSheet.Rows.Item(2).RowHeight = 43.5

# B2: the estimate label, word-wrapped
$newSheet.Range("B2").Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$newSheet.Range("B2").WrapText = $true

# C2: the estimate total
$newSheet.Range("C2").Value = 314

# Page setup for the new sheet
$newSheet.PageSetup.Orientation = 1

# Make "Estimates" the active/selected sheet and cell, like in the saved file
$newSheet.Activate()
$newSheet.Range("C2").Select()
